$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column B (Total) updates
$ws.Range("B3").Value = 6543.656484598552
$ws.Range("B4").Value = 3019.34496453037
$ws.Range("B5").Value = 2058.906273972604
$ws.Range("B6").Value = 3768.334294520545
$ws.Range("B7").Value = 6107.436147945209
$ws.Range("B8").Value = 8212.240854109599
$ws.Range("B9").Value = 12582.76608219179

# Column D (Community) updates
$ws.Range("D3").Value = 270.6623123287663
$ws.Range("D4").Value = 94.56333487926334
$ws.Range("D5").Value = 76.18662739726032
$ws.Range("D6").Value = 188.68975479452
$ws.Range("D7").Value = 271.2249986301366
$ws.Range("D8").Value = 356.041101369863
$ws.Range("D9").Value = 432.6294383561631

# Totals / ratios
$ws.Range("F10").Value = 9422865.337821921
$ws.Range("G11").Value = 0.7109339862006124
$ws.Range("F12").Value = 389753.7297534242
$ws.Range("G12").Value = 0.04136254905278261
$ws.Range("G13").Value = 0.247703464746605
